# Update imputed values in the KNN result data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value  = -13.463
$ws.Range("A8").Value  = -21.255
$ws.Range("A10").Value = -21.047
$ws.Range("A12").Value = -21.649
$ws.Range("B13").Value = 6.753
$ws.Range("A18").Value = -21.649
$ws.Range("C20").Value = -13.041
$ws.Range("A25").Value = -21.534
